$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.278945326805115
$ws.Range("B1").Value = 2.322174072265625
$ws.Range("C1").Value = 3.89882755279541
$ws.Range("D1").Value = 2.886640071868896
$ws.Range("E1").Value = 1.340027809143066
